# Overall_Running_Metadata_for_All_LCMSMS_Jobs.xlsx
# - "Job to Run" sheet: change the active job's name from
#   "Anid_HE_TJGIp11_pos_2018" to "Anid_HE_TJGIp4_TMM_pos_2018"
# - "All" sheet: append a row for the new job (copy of the EXP/CTRL
#   replicate counts and ionization of the previous job), drop the
#   leftover blank formatted rows, and clean up row 2's one-off
#   row height/format override.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Job to Run")
$ws2 = $wb.Worksheets.Item("All")

# ---------------------------------------------------------------
# "Job to Run": update the job name used for this run
# ---------------------------------------------------------------
$ws1.Range("A2").Value = "Anid_HE_TJGIp4_TMM_pos_2018"

# widen column A (job name is much longer now) and nudge column B
$ws1.Columns.Item(1).ColumnWidth = 65.7
$ws1.Columns.Item(2).ColumnWidth = 15.7

$ws1.Activate() | Out-Null
$ws1.Range("A2:D2").Select() | Out-Null

# ---------------------------------------------------------------
# "All": drop old trailing placeholder rows (3:6) and record the
# new job as a fresh row, keeping the same layout/styling as row 2
# ---------------------------------------------------------------
$ws2.Range("A3:D6").EntireRow.Delete() | Out-Null

$row2 = $ws2.Rows.Item(2)
$row2.ClearFormats() | Out-Null
$ws2.Range("A2:D2").VerticalAlignment = -4108
$row2.AutoFit() | Out-Null

$ws2.Range("A3").Value = "Anid_HE_TJGIp4_TMM_pos_2018"
$ws2.Range("B3").Value = 3
$ws2.Range("C3").Value = 3
$ws2.Range("D3").Value = "POS"
$ws2.Range("A3:D3").VerticalAlignment = -4108

$ws2.Columns.Item(1).ColumnWidth = 25.02
$ws2.Columns.Item(2).ColumnWidth = 15.7

$ws2.Activate() | Out-Null
$ws2.Range("B12").Select() | Out-Null

# restore "Job to Run" as the visible/active sheet
$ws1.Activate() | Out-Null
